$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update L3 value: 2020 -> 2021
$ws.Range("L3").Value = 2021

# Add new column M with data for row 3 and row 4, copying style from column L
$ws.Range("L3").Copy() | Out-Null
$ws.Range("M3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("M3").Value = 2022

$ws.Range("L4").Copy() | Out-Null
$ws.Range("M4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("M4").Value = 6.18

$excel.CutCopyMode = 0

# Update the selection to match the new active cell
$ws.Range("M9").Select() | Out-Null
